$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values look like plain numbers (e.g. "316.23") and Excel would
# auto-convert them to doubles (losing exact text / trailing zeros) unless
# we force the cell to Text format first, then restore its original style.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '44.658.95'
$ws.Range('E2').Value = '  +3.78%  '
Set-TextValue $ws.Range('D3') '2.422.35'
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue $ws.Range('D5') '316.23'
$ws.Range('E5').Value = '  +4.14%  '
Set-TextValue $ws.Range('D6') '101.63'
$ws.Range('E6').Value = '  +6.26%  '
Set-TextValue $ws.Range('D7') '0.513'
$ws.Range('E7').Value = '  +2.15%  '
$ws.Range('E8').Value = '  -0.08%  '
Set-TextValue $ws.Range('D9') '0.525'
$ws.Range('E9').Value = '  +8.97%  '
Set-TextValue $ws.Range('D10') '35.46'
$ws.Range('E10').Value = '  +3.48%  '
$ws.Range('E11').Value = '  +1.72%  '
Set-TextValue $ws.Range('D12') '18.88'
$ws.Range('E12').Value = '  +1.97%  '
$ws.Range('E13').Value = '  -2.72%  '
Set-TextValue $ws.Range('D14') '6.95'
$ws.Range('E14').Value = '  +2.98%  '
Set-TextValue $ws.Range('D15') '2.799.07'
Set-TextValue $ws.Range('D16') '2.399.30'
$ws.Range('E16').Value = '  +2.28%  '
Set-TextValue $ws.Range('D17') '0.833'
$ws.Range('E17').Value = '  +4.42%  '
Set-TextValue $ws.Range('D18') '44.493.85'
$ws.Range('E18').Value = '  +3.51%  '
Set-TextValue $ws.Range('D19') '12.34'
$ws.Range('E19').Value = '  +3.28%  '
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('E21').Value = '  +3.80%  '
Set-TextValue $ws.Range('D22') '68.78'
$ws.Range('E22').Value = '  +0.95%  '
Set-TextValue $ws.Range('D23') '241.90'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('E24').Value = '  +4.71%  '
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('E26').Value = '  -0.09%  '
Set-TextValue $ws.Range('D27') '25.23'
$ws.Range('E27').Value = '  +2.69%  '
Set-TextValue $ws.Range('D28') '2.29'
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('E29').Value = '  +1.48%  '
Set-TextValue $ws.Range('D30') '33.50'
$ws.Range('E30').Value = '  +4.45%  '
Set-TextValue $ws.Range('D31') '48.51'
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('E32').Value = '  +19.20%  '
Set-TextValue $ws.Range('D33') '19.50'
$ws.Range('E33').Value = '  +10.92%  '
$ws.Range('E34').Value = '  +3.38%  '
Set-TextValue $ws.Range('D35') '0.0777'
$ws.Range('E35').Value = '  +8.42%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('E38').Value = '  +3.09%  '
Set-TextValue $ws.Range('D39') '2.88'
$ws.Range('E39').Value = '  +1.10%  '
Set-TextValue $ws.Range('D40') '121.23'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('E42').Value = '  -2.71%  '
Set-TextValue $ws.Range('D43') '20.95'
$ws.Range('E43').Value = '  -1.16%  '
$ws.Range('E44').Value = '  +4.20%  '
Set-TextValue $ws.Range('D45') '1.944.96'
$ws.Range('E45').Value = '  +0.54%  '
Set-TextValue $ws.Range('D47') '2.95'
$ws.Range('E47').Value = '  +8.44%  '
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  +11.09%  '
Set-TextValue $ws.Range('D50') '54.27'
$ws.Range('E50').Value = '  +5.92%  '
Set-TextValue $ws.Range('D51') '75.07'
$ws.Range('E51').Value = '  +4.68%  '

Write-Host "Updated cryptos list"
